# B_survival_07172024.xlsx — fill in the "surveyor" column (G) for the
# 6/20/24 bag-emptying rows (32-41) that were left blank. Odd-numbered
# bags were counted by Erik, even-numbered ones by Grace (matches the
# alternating pattern already used earlier in the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

$surveyorByRow = @{
    32 = "Erik"
    33 = "Erik"
    34 = "Grace"
    35 = "Grace"
    36 = "Grace"
    37 = "Erik"
    38 = "Erik"
    39 = "Grace"
    40 = "Grace"
    41 = "Grace"
}

foreach ($row in 32..41) {
    $ws.Cells.Item($row, 7).Value = $surveyorByRow[$row]
}

# Scroll the view down and leave the selection where it ended up, same as
# the saved workbook (best effort — the engine may not persist scroll
# position outside of frozen panes).
$ws.Range("A16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("K25").Select() | Out-Null
